$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: every "Ready for handoff" cell becomes "Handed back: in
#    sync with en-US" now that the handback round-trip completed.
# ---------------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Widen the "Status"/"Latest Handback File" style columns that now need
#    to fit the longer strings.
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("J1").ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. Record the handback results for the "a.md" row (row 2 / row 3 are the
#    same source file a.md) on each language sheet: the target file that
#    was handed back, the handback xliff file, and (for de-de, the last
#    language processed) the handback timestamp.
# ---------------------------------------------------------------------------
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0648b9a1cb95065baecf16f02f54dc7abbce9102/e2e/a.md"

# zh-cn: Latest Target File (I) + Latest Handback File (J)
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")

# de-de: Latest Target File (I) + Latest Handback File (J) + Latest Handback DateTime (K)
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsDeDe.Range("K2").Value = "2016-08-19 14:43:37"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$wsDeDe.Range("K3").Value = "2016-08-19 14:43:37"

# zh-cn Latest Handback DateTime (K) was the first language handed back.
$wsZhCn.Range("K2").Value = "2016-08-19 14:43:31"
$wsZhCn.Range("K3").Value = "2016-08-19 14:43:31"
